# verzia z dna 14.1.25 = obhajoba
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# APR_24 sheet: update an expense value (B7) which feeds the F2 "stav na
# uctu" formula (SUM(B:B)-SUM(A:A)), and move the selection to B8.
# ---------------------------------------------------------------------------
$apr = $wb.Worksheets.Item("APR_24")
$apr.Range("B7").Value = 800
$apr.Range("B8").Select() | Out-Null

# ---------------------------------------------------------------------------
# JAN_25 sheet: fix the mislabeled header, update amounts/dates for the
# existing rows, and append a new expense row (row 8).
#
# The order in which *new* text values are written matters because new
# entries are appended to the shared-string table in first-use order; we
# mirror the order the author must have typed them in (descriptions for
# rows 3-7 first, then the dates for rows 2-7, then the new row's data).
# ---------------------------------------------------------------------------
$jan25 = $wb.Worksheets.Item("JAN_25")

# Fix header A1 ("AP" -> "Vydaj", to match the other month sheets)
$jan25.Range("A1").Value = "Vydaj"

# Updated amounts
$jan25.Range("A2").Value = 500
$jan25.Range("A3").Value = 2000
$jan25.Range("A4").Value = 50
$jan25.Range("B5").Value = 1500
$jan25.Range("B6").Value = 1700
$jan25.Range("B7").Value = 300

# Updated descriptions (column D) for rows 3-7
$jan25.Range("D3").Value = "rozbite dvere"
$jan25.Range("D4").Value = "benzin"
$jan25.Range("D5").Value = "zbierka"
$jan25.Range("D6").Value = "pozicka"
$jan25.Range("D7").Value = "poistenie dveri"

# Updated dates (column C) for rows 2-7. Some of these ("1.1.2025",
# "5.1.2025", "7.1.2025", "10.1.2025") look like valid M.D.Y dates to
# Excel's parser and would silently be converted to date serials, so force
# text entry for those with a temporary "@" number format, then restore the
# original (General) cell style afterwards by pasting the formatting from a
# cell that was left untouched.
$jan25.Range("C2").NumberFormat = "@"
$jan25.Range("C2").Value = "1.1.2025"
$jan25.Range("C3").NumberFormat = "@"
$jan25.Range("C3").Value = "5.1.2025"
$jan25.Range("C4").NumberFormat = "@"
$jan25.Range("C4").Value = "7.1.2025"
$jan25.Range("C5").NumberFormat = "@"
$jan25.Range("C5").Value = "10.1.2025"
$jan25.Range("C6").Value = "16.1.2025"
$jan25.Range("C7").Value = "26.1.2025"

$jan25.Range("C7").Copy() | Out-Null
$jan25.Range("C2:C5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# New row 8: copy formatting from the row-4 pattern (expense row), then set
# its values - date and description last, matching the append order of the
# new shared strings.
$jan25.Range("A4:D4").Copy($jan25.Range("A8:D8")) | Out-Null
$jan25.Range("A8").Value = 650
$jan25.Range("B8").Value = "—"
$jan25.Range("C8").Value = "30.1.2025"
$jan25.Range("D8").Value = "kupa kosacky"

# ---------------------------------------------------------------------------
# Sheet1 (summary): the Jan_2025 row pointed at the wrong month (JAN_24)
# by mistake - repoint it at JAN_25.
# ---------------------------------------------------------------------------
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Range("B14").Formula = "=JAN_25!F2"
$sheet1.Range("B15").Select() | Out-Null

# ---------------------------------------------------------------------------
# View state: JAN_25 is the tab that should end up active/selected.
# ---------------------------------------------------------------------------
$jan25.Activate()
$jan25.Range("H12").Select() | Out-Null
